# Kilroy Equipment List 2017 - add USB/Axis camera rows, front motor
# controllers, ring light relay, PDB and aligning-button entries; bump
# the PWM port used by the rear-right motor controller row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Speed Controllers section ---
$ws.Range("E4").Value = "PWM 4"

$ws.Range("C6").Value = "Right Front Motor Controller"
$ws.Range("D6").Value = "rightFrontMotor"
$ws.Range("E6").Value = "PWM 2"

$ws.Range("C7").Value = "Left Front Motor Controller"
$ws.Range("D7").Value = "leftFrontMotor"
$ws.Range("E7").Value = "PWM 1"

# --- Relays section ---
$ws.Range("C44").Value = "Ring Light Relay"
$ws.Range("D44").Value = "ringlightRelay"
$ws.Range("E44").Value = "0"

# --- PDB (Power Distribution Board) section ---
$ws.Range("C62").Value = "Power Distribution Board/Panel"
$ws.Range("D62").Value = "pdp"
$ws.Range("E62").Value = "0"

# --- Operator Controls section ---
$ws.Range("C72").Value = "Set Rotation Value"
$ws.Range("E72").Value = "Right Driver Trigger"

$ws.Range("C73").Value = "Stop Aligning/ set isAligning to False"
$ws.Range("E73").Value = "Left Operator 7"

$ws.Range("C74").Value = "Start Aligning/ set isAligning to True"
$ws.Range("E74").Value = "Left Operator 8"

$ws.Range("C75").Value = "Take Single Picture "
$ws.Range("E75").Value = "Left Operator 8"

# --- Miscellaneous section: Axis + USB cameras ---
$ws.Range("C109").Value = "Axis Camera"
$ws.Range("D109").Value = "axisCamera"
$ws.Range("E109").Value = "10.3.39.11"

$ws.Range("C110").Value = "USB Camera 0"
$ws.Range("D110").Value = "cam0"
$ws.Range("E110").Value = "n/a"

$ws.Range("C111").Value = "USB Camera 1"
$ws.Range("D111").Value = "cam1"
$ws.Range("E111").Value = "n/a"

# Restore the view: scroll the frozen pane down to the newly-added rows
# and leave the selection on the last edited cell.
$ws.Range("A104").Select() | Out-Null
$ws.Range("E111").Select() | Out-Null
